# Apply the "Error Calculations and Plots" data corrections to the
# missing_data worksheet: a couple of measurement rows are dropped
# entirely (their data shifts up one row), and a number of individual
# cells flip between having a numeric reading and being blank
# (missing/imputed), matching the committed OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 ("RM 232") is removed outright - every row below shifts up by one.
$ws.Rows(26).Delete()
# The row that is now at 27 ("SC 92", originally row 28) is removed too -
# remaining rows shift up again, so the sheet ends at row 33 (A1:F33).
$ws.Rows(27).Delete()

# --- Cell-level value/blank changes in the untouched-row-number block (rows 2-25) ---
$ws.Cells.Item(2,3).Value = 14.9      # RM 2   - B -> 14.9
$ws.Cells.Item(3,4).Value = ""        # RM 8   - C -> blank
$ws.Cells.Item(4,4).Value = -15.4     # RM 9   - C -> -15.4
$ws.Cells.Item(5,4).Value = ""        # RM 14  - C -> blank
$ws.Cells.Item(6,3).Value = ""        # RM 21  - B -> blank
$ws.Cells.Item(8,4).Value = ""        # RM 38  - C -> blank
$ws.Cells.Item(12,3).Value = 12.5     # RM 81  - B -> 12.5
$ws.Cells.Item(14,3).Value = ""       # RM 90  - B -> blank
$ws.Cells.Item(15,4).Value = -15.2    # RM 95  - C -> -15.2
$ws.Cells.Item(18,4).Value = -15.2    # RM 120 - C -> -15.2
$ws.Cells.Item(19,4).Value = ""       # RM 125 - C -> blank
$ws.Cells.Item(20,3).Value = 12.5     # RM 134 - B -> 12.5
$ws.Cells.Item(21,3).Value = 12.7     # RM 135 - B -> 12.7
$ws.Cells.Item(22,4).Value = ""       # RM 138 - C -> blank
$ws.Cells.Item(23,3).Value = ""       # RM 140 - B -> blank
$ws.Cells.Item(23,4).Value = -13.9    # RM 140 - C -> -13.9
$ws.Cells.Item(24,3).Value = ""       # RM 142a- B -> blank
$ws.Cells.Item(25,4).Value = -15.5    # RM 145 - C -> -15.5

# --- Cell-level value/blank changes in the shifted block (now rows 26-33) ---
$ws.Cells.Item(26,2).Value = -20.2    # SC 5   - A -> -20.2
$ws.Cells.Item(27,2).Value = ""       # SC 101 - A -> blank
$ws.Cells.Item(27,4).Value = ""       # SC 101 - C -> blank
$ws.Cells.Item(30,2).Value = -19.7    # SC 120 - A -> -19.7
$ws.Cells.Item(31,3).Value = 15.3     # SC 132 - B -> 15.3
$ws.Cells.Item(32,2).Value = ""       # SC 193 - A -> blank
$ws.Cells.Item(33,3).Value = 10.4     # SC 232 - B -> 10.4
